$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the card number / client for the last two kept rows.
# The leading "'" forces the numeric-looking card numbers to stay text
# (matching the rest of column A), and resetting the style back to
# "Normal" afterwards avoids leaving a stray number-format override on
# the cell.
$ws.Range("A22").Value = "'364468"
$ws.Range("A22").Style = "Normal"
$ws.Range("B22").Value = "Kosuke Yokono"

$ws.Range("A23").Value = "'584106"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = "Hana Murata"

# Remove the remaining now-obsolete client rows (24-31).
$ws.Rows("24:31").Delete()
